$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Swap B14 and B15 values
$ws.Range("B14").Value = "供销大集"
$ws.Range("B15").Value = "国轩高科"

# Rotate C15, C16, C17 values
$ws.Range("C15").Value = "领益智造"
$ws.Range("C16").Value = "上海电力"
$ws.Range("C17").Value = "天普股份"
